$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 319
$ws.Range('F4').Value = 3274
$ws.Range('F7').Value = 333
$ws.Range('F8').Value = 7651
$ws.Range('F11').Value = 23
$ws.Range('F12').Value = 41
$ws.Range('F14').Value = 671
$ws.Range('F15').Value = 1110
$ws.Range('F16').Value = 1047
$ws.Range('F17').Value = 1012
$ws.Range('F18').Value = 165
$ws.Range('F19').Value = 1625
$ws.Range('F21').Value = 6056
$ws.Range('F24').Value = 1002
$ws.Range('F27').Value = 4201
$ws.Range('F28').Value = 3845
$ws.Range('F29').Value = 292
$ws.Range('F30').Value = 97
$ws.Range('F32').Value = 269
$ws.Range('F33').Value = 1028
$ws.Range('F34').Value = 1029
$ws.Range('F35').Value = 1016
$ws.Range('F36').Value = 83
$ws.Range('F38').Value = 416
$ws.Range('F42').Value = 383
$ws.Range('F43').Value = 315
$ws.Range('F44').Value = 1098
$ws.Range('F45').Value = 465
$ws.Range('F46').Value = 3104
$ws.Range('F47').Value = 74
$ws.Range('F48').Value = 349
$ws.Range('F49').Value = 5

# Sheet: 演出
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F5').Value = 3
$ws.Range('F6').Value = 116
$ws.Range('F9').Value = 375
$ws.Range('F10').Value = 606
$ws.Range('F15').Value = 245
$ws.Range('F17').Value = 1
$ws.Range('F20').Value = 152
$ws.Range('F28').Value = 6126
$ws.Range('G28').Value = 980
$ws.Range('B29').Value = '2024-10-18'
$ws.Range('C29').Value = '上海·ROOKiEZ is PUNK`D 「Reignite Youth （重燃青春）」2024 CHINA Tour '
$ws.Range('D29').Value = '虹许路731号4号楼 THE BOXX•城市乐园'
$ws.Range('E29').Value = '2024.10.18 20:30-10.18 22:00'
$ws.Range('F29').Value = 44
$ws.Range('G29').Value = 259
$ws.Range('H29').Value = 'https://show.bilibili.com/platform/detail.html?id=91376'
$ws.Range('I29').Value = '//i1.hdslb.com/bfs/openplatform/202408/pZdI02BJ1724735899119.jpeg'
$ws.Range('B30').Value = '2024-10-26'
$ws.Range('C30').Value = '上海·【早鸟4折】“海上钢琴师”一生必听经典电影主题音乐会'
$ws.Range('D30').Value = '南京西路1376号 上海商城剧院'
$ws.Range('E30').Value = '2024.10.26 15:00-10.26 16:30'
$ws.Range('F30').Value = 2
$ws.Range('G30').Value = 48
$ws.Range('H30').Value = 'https://show.bilibili.com/platform/detail.html?id=91375'
$ws.Range('I30').Value = '//i1.hdslb.com/bfs/openplatform/202408/qfPgppOK1724743485013.jpeg'
$ws.Range('B31').Value = '2024-10-31'
$ws.Range('C31').Value = '上海·苏菲•珊曼妮2024巡回演唱会'
$ws.Range('D31').Value = '重庆南路308号3楼 上海MAO LIVEHOUSE'
$ws.Range('E31').Value = '2024.10.31 20:00-10.31 21:40'
$ws.Range('F31').Value = 8
$ws.Range('G31').Value = 380
$ws.Range('H31').Value = 'https://show.bilibili.com/platform/detail.html?id=87918'
$ws.Range('I31').Value = '//i2.hdslb.com/bfs/openplatform/202406/RhhjOqDY1718160939240.jpeg'
$ws.Range('B32').Value = '2024-11-02'
$ws.Range('C32').Value = '上海·欢迎来到绵羊咖啡屋! 中国第二回 仲村宗悟 梶原岳人'
$ws.Range('D32').Value = '茂名南路57号近长乐路 上海兰心大戏院'
$ws.Range('E32').Value = '2024.11.02 12:00-11.02 21:30'
$ws.Range('F32').Value = 52
$ws.Range('H32').Value = 'https://show.bilibili.com/platform/detail.html?id=91176'
$ws.Range('I32').Value = '//i0.hdslb.com/bfs/openplatform/202408/QtV99HO81724384169942.jpeg'
$ws.Range('B33').Value = '2024-11-15'
$ws.Range('C33').Value = '上海·“法国姐姐”乔伊丝·乔纳森《小意思》'
$ws.Range('D33').Value = '高青西路777号 上海前滩31演艺中心'
$ws.Range('E33').Value = '2024.11.15 19:30-11.15 21:00'
$ws.Range('F33').Value = 0
$ws.Range('G33').Value = 280
$ws.Range('H33').Value = 'https://show.bilibili.com/platform/detail.html?id=91619'
$ws.Range('I33').Value = '//i1.hdslb.com/bfs/openplatform/202408/VnZEk71H1725014748758.jpeg'

# Sheet: 本地生活
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F2').Value = 122
$ws.Range('F6').Value = 1952
$ws.Range('F9').Value = 1205
$ws.Range('F12').Value = 524
$ws.Range('F13').Value = 2061
$ws.Range('F14').Value = 8767
$ws.Range('F15').Value = 904
$ws.Range('F16').Value = 57

# Sheet: 全部类型
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 319
$ws.Range('F4').Value = 1952
$ws.Range('F5').Value = 333
$ws.Range('F7').Value = 7651
$ws.Range('F8').Value = 1205
$ws.Range('F10').Value = 524
$ws.Range('F12').Value = 23
$ws.Range('F13').Value = 41
$ws.Range('F15').Value = 904
$ws.Range('F16').Value = 606
$ws.Range('F17').Value = 57
$ws.Range('F18').Value = 671
$ws.Range('F19').Value = 1110
$ws.Range('F20').Value = 1047
$ws.Range('F21').Value = 1012
$ws.Range('F23').Value = 165
$ws.Range('F24').Value = 245
$ws.Range('F26').Value = 6056
$ws.Range('F28').Value = 1002
$ws.Range('F30').Value = 4201
$ws.Range('F31').Value = 292
$ws.Range('F32').Value = 97
$ws.Range('F34').Value = 269
$ws.Range('F35').Value = 1028
$ws.Range('F36').Value = 1029
$ws.Range('F37').Value = 1016
$ws.Range('F40').Value = 152
$ws.Range('F43').Value = 383
$ws.Range('F44').Value = 315
$ws.Range('F46').Value = 465
$ws.Range('F47').Value = 3104
$ws.Range('F48').Value = 74
$ws.Range('F49').Value = 6126
$ws.Range('G49').Value = 980
